# Daily attendance processing - 2026-01-16 09:42:15
# Reverses the order of the comma-separated "Recorded By" entries (column G)
# for every row where the recorded-by list contains "System" (any casing).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        if ($value -like "*,*" -and $value -like "*ystem*") {
            $parts = $value -split ","
            $trimmed = @()
            foreach ($p in $parts) {
                $trimmed += $p.Trim()
            }

            $reversed = @()
            for ($i = $trimmed.Count - 1; $i -ge 0; $i--) {
                $reversed += $trimmed[$i]
            }

            $newValue = [string]::Join(", ", $reversed)
            $cell.Value2 = $newValue
        }
    }
}
